$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new sign-up row (row 9): name, email (hyperlinked), phone ---
$ws.Range("A9").Value = "Afonso Hermenegildo"
$ws.Range("B9").Value = "afonsosousah@tec2med.com"
$ws.Range("C9").Value = 912345678

# Give the new name/phone cells the same style used by the rows above them
$ws.Range("A9").Style = $ws.Range("A8").Style

# Turn the new email cell into a mailto hyperlink, matching the existing
# hyperlinked email cell above it (B8)
$ws.Hyperlinks.Add($ws.Range("B9"), "mailto:afonsosousah@tec2med.com")

# --- Row heights / column width refresh (matches the re-saved layout) ---
$ws.Range("A3:C9").RowHeight = 17.4
$ws.Range("B:B").ColumnWidth = 34.46

# --- Selection moved to D9 ---
$ws.Range("D9").Select()
